# This script applies the odds updates described in the commit diff
# ("Atualizando o arquivo XLSX") to the Betfair Back/Lay odds sheet for 2025-11-12.
# It updates the numeric odds/stake values across rows 2-6 to their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)


# Row 2 updates
$ws.Range("J2").Value = 1.03
$ws.Range("N2").Value = 1.25
$ws.Range("P2").Value = 1.24
$ws.Range("R2").Value = 1.18
$ws.Range("T2").Value = 1.04
$ws.Range("U2").Value = 1.04

# Row 3 updates
$ws.Range("F3").Value = 2.22
$ws.Range("G3").Value = 2.6
$ws.Range("H3").Value = 3.35
$ws.Range("I3").Value = 4.2
$ws.Range("J3").Value = 2.94
$ws.Range("K3").Value = 3.45
$ws.Range("L3").Value = 1.5
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 2.5
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 1.58
$ws.Range("Q3").Value = 2.32
$ws.Range("R3").Value = 1.21
$ws.Range("T3").Value = 1.94
$ws.Range("U3").Value = 1.83
$ws.Range("V3").Value = 1.32
$ws.Range("W3").Value = 1.71
$ws.Range("X3").Value = 11
$ws.Range("Y3").Value = 12.5
$ws.Range("Z3").Value = 27
$ws.Range("AB3").Value = 9.6
$ws.Range("AC3").Value = 8.6
$ws.Range("AD3").Value = 18
$ws.Range("AE3").Value = 60
$ws.Range("AF3").Value = 16.5
$ws.Range("AG3").Value = 13.5
$ws.Range("AH3").Value = 25
$ws.Range("AJ3").Value = 42
$ws.Range("AK3").Value = 38
$ws.Range("AN3").Value = 38
$ws.Range("AO3").Value = 75

# Row 4 updates
$ws.Range("G4").Value = 2.08
$ws.Range("H4").Value = 4
$ws.Range("J4").Value = 3.25
$ws.Range("L4").Value = 1.47
$ws.Range("N4").Value = 2.64
$ws.Range("O4").Value = 1.41
$ws.Range("P4").Value = 1.65
$ws.Range("Q4").Value = 2.22
$ws.Range("R4").Value = 1.24
$ws.Range("S4").Value = 3.8
$ws.Range("T4").Value = 1.98
$ws.Range("U4").Value = 1.81
$ws.Range("V4").Value = 1.22
$ws.Range("W4").Value = 1.92
$ws.Range("AC4").Value = 9.4
$ws.Range("AG4").Value = 13

# Row 5 updates
$ws.Range("F5").Value = 1.71
$ws.Range("G5").Value = 1.73
$ws.Range("I5").Value = 6.4
$ws.Range("K5").Value = 3.95
$ws.Range("N5").Value = 3.5
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 1.86
$ws.Range("Q5").Value = 2.1
$ws.Range("R5").Value = 1.32
$ws.Range("T5").Value = 2.06
$ws.Range("W5").Value = 2.36
$ws.Range("Y5").Value = 18.5
$ws.Range("AC5").Value = 8.800000000000001

# Row 6 updates
$ws.Range("F6").Value = 5.2
$ws.Range("G6").Value = 6
$ws.Range("H6").Value = 1.77
$ws.Range("J6").Value = 3.4
$ws.Range("K6").Value = 3.75
$ws.Range("L6").Value = 1.39
$ws.Range("N6").Value = 3.1
$ws.Range("O6").Value = 1.4
$ws.Range("P6").Value = 1.72
$ws.Range("Q6").Value = 2.08
$ws.Range("R6").Value = 1.26
$ws.Range("S6").Value = 4
$ws.Range("T6").Value = 1.98
$ws.Range("U6").Value = 1.81
$ws.Range("W6").Value = 1.2
$ws.Range("AI6").Value = 60
$ws.Range("AK6").Value = 110
$ws.Range("AL6").Value = 120
